# Update "想去人数" (number of people interested) values in the
# "展览" and "全部类型" sheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 3143
    6  = 1761
    8  = 102
    9  = 38
    11 = 1454
    13 = 576
    15 = 86
    21 = 98
    23 = 3432
    27 = 73
    28 = 21
    30 = 1164
    31 = 126
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
